$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Insert a new row at position 10, shifting existing rows (10-21) down to (11-22)
$ws.Rows.Item(10).Insert()

# Fill in the new row's data: SyndName / 聯貸名稱 / NVARCHAR2 / length 60
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "SyndName"
$ws.Range("C10").Value = "聯貸名稱"
$ws.Range("D10").Value = "NVARCHAR2"
$ws.Range("E10").Value = 60

# Copy formatting from neighboring rows that already have the right look
$ws.Range("B11").Copy()
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)

$ws.Range("D9:E9").Copy()
$ws.Range("D10:E10").PasteSpecial(-4122)

$ws.Range("F11").Copy()
$ws.Range("F10").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Activate()
$ws.Range("G10").Select()
